$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 33) with the latest monthly data point (01-08-2021).
$row = 33

# Column A ("Serie") holds a label that looks like a date, e.g. "01-08-2021".
# Assigning that text straight to a cell's Value/Formula makes Excel
# auto-recognize it as a real date (storing a serial number and tagging the
# cell with a date number-format style). To keep it as plain text - exactly
# like the existing cells in that column - build it first as a text formula
# result (which Excel keeps as a string) on a scratch cell, then copy/paste
# only the resulting value into the target cell. That avoids both the
# unwanted date conversion and leaves no extra cell style behind.
$helper = $ws.Cells.Item(100, 1)
$helper.Formula = '="01-08-2021"'
$helper.Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$helper.Clear()

$ws.Cells.Item($row, 2).Value = 84660
$ws.Cells.Item($row, 3).Value = 15861
$ws.Cells.Item($row, 4).Value = 8832
$ws.Cells.Item($row, 5).Value = 5538
$ws.Cells.Item($row, 6).Value = 6199
$ws.Cells.Item($row, 7).Value = 6549
$ws.Cells.Item($row, 8).Value = 17996
$ws.Cells.Item($row, 9).Value = 14438
$ws.Cells.Item($row, 10).Value = 9247
